# Update "想去人数" (F column) figures for the same set of events on both
# the "展览" sheet and the "全部类型" sheet, matching the refreshed data
# output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# 展览 sheet: rows 2,3,4,7,17,22,28,34
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14852
$ws1.Range("F3").Value = 18436
$ws1.Range("F4").Value = 146
$ws1.Range("F7").Value = 221
$ws1.Range("F17").Value = 1404
$ws1.Range("F22").Value = 7637
$ws1.Range("F28").Value = 5944
$ws1.Range("F34").Value = 5280

# 全部类型 sheet: same events, but shifted by extra rows, so row numbers differ
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14852
$ws4.Range("F3").Value = 18436
$ws4.Range("F4").Value = 146
$ws4.Range("F7").Value = 221
$ws4.Range("F17").Value = 1404
$ws4.Range("F23").Value = 7637
$ws4.Range("F30").Value = 5944
$ws4.Range("F36").Value = 5280
